$wb = $excel.ActiveWorkbook

# --- 1. Rename the existing "mars" sheet to "March" ---
$march = $wb.Worksheets.Item(1)
$march.Name = "March"

# --- 2. Update row 6 on "March": was Clothing/genser/.../Checkings,
#        now Transportation/atb/.../Savings ---
# Force text (shared-string) storage for numeric/date-looking values so the
# cells keep the same "t=s" shared-string representation as the rest of the
# sheet, then strip the temporary number-format style back off again.
$marchRange = $march.Range("A6:E7")
$marchRange.NumberFormat = "@"

$march.Range("A6").Value = "Transportation"
$march.Range("B6").Value = "atb"
$march.Range("C6").Value = "2023-03-20"
$march.Range("D6").Value = "1000.0"
$march.Range("E6").Value = "Savings"

# --- 3. Add new row 7 on "March" with the original Clothing/genser entry ---
$march.Range("A7").Value = "Clothing"
$march.Range("B7").Value = "genser"
$march.Range("C7").Value = "2023-03-20"
$march.Range("D7").Value = "1000.0"
$march.Range("E7").Value = "Checkings"

$marchRange.ClearFormats()

# --- 4. Add the new "April" sheet right after "March" ---
$april = $wb.Worksheets.Add($null, $march)
$april.Name = "April"

$aprilRange = $april.Range("A1:E2")
$aprilRange.NumberFormat = "@"

$april.Range("A1").Value = "Category"
$april.Range("B1").Value = "Name"
$april.Range("C1").Value = "Date"
$april.Range("D1").Value = "Price"
$april.Range("E1").Value = "Account"

$april.Range("A2").Value = "Entertainment"
$april.Range("B2").Value = "playstation"
$april.Range("C2").Value = "2023-04-12"
$april.Range("D2").Value = "2000.0"
$april.Range("E2").Value = "Savings"

$aprilRange.ClearFormats()

# --- 5. Keep "March" as the active/selected tab (activeTab stays 0) ---
$march.Activate()
